$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 183.33333
$ws.Range("I4").Value = 183.33333
$ws.Range("K4").Value = 183.33333
$ws.Range("M4").Value = -69.33332999999999

$ws.Range("H43").Value = 721.8889
$ws.Range("I43").Value = 679.4
$ws.Range("J43").Value = 775
$ws.Range("K43").Value = 679.4
$ws.Range("L43").Value = 775
$ws.Range("M43").Value = -610.4
$ws.Range("N43").Value = -913

$ws.Range("H98").Value = 705.4167
$ws.Range("I98").Value = 705.4167
$ws.Range("K98").Value = 705.4167
$ws.Range("M98").Value = 792.5833

$ws.Range("M112").ClearContents()
$ws.Range("H112").Value = 3099.2856
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3099.2856
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9297.856800000001
$ws.Range("N112").Value = -11513.8568

$ws.Range("H122").Value = 705.4167
$ws.Range("I122").Value = 705.4167
$ws.Range("K122").Value = 2116.2501
$ws.Range("M122").Value = 333.7498999999998

$ws.Range("M127").ClearContents()
$ws.Range("H127").Value = 1687.7142
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1687.7142
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 5063.142599999999
$ws.Range("N127").Value = -14983.1426

$ws.Range("H129").Value = 838.1622
$ws.Range("J129").Value = 879.55225
$ws.Range("L129").Value = 2638.65675
$ws.Range("N129").Value = -12638.65675

$ws.Range("H137").Value = 1981.1017
$ws.Range("I137").Value = 2158.121
$ws.Range("J137").Value = 1756.4231
$ws.Range("K137").Value = 6474.363
$ws.Range("L137").Value = 5269.2693
$ws.Range("M137").Value = -3924.363
$ws.Range("N137").Value = -10369.2693

$ws.Range("H138").Value = 2622
$ws.Range("I138").Value = 899.6667
$ws.Range("J138").Value = 3000.0732
$ws.Range("K138").Value = 2699.0001
$ws.Range("L138").Value = 9000.2196
$ws.Range("M138").Value = 2440.9999
$ws.Range("N138").Value = -19280.2196

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2320.5076
$ws.Range("I32").Value = 1765.4839
$ws.Range("K32").Value = 1765.4839
$ws.Range("M32").Value = -1478.4839

$ws.Range("H45").Value = 3204.92
$ws.Range("I45").Value = 2388.625
$ws.Range("K45").Value = 2388.625
$ws.Range("M45").Value = -2011.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2408.1667
$ws.Range("I20").Value = 2546.2144
$ws.Range("J20").Value = 1925
$ws.Range("K20").Value = 2546.2144
$ws.Range("L20").Value = 1925
$ws.Range("M20").Value = -2299.2144
$ws.Range("N20").Value = -2419

$ws.Range("H82").Value = 15200.7
$ws.Range("I82").Value = 5184
$ws.Range("J82").Value = 55267.5
$ws.Range("K82").Value = 5184
$ws.Range("L82").Value = 55267.5
$ws.Range("M82").Value = -4801
$ws.Range("N82").Value = -56033.5

$ws.Range("H85").Value = 15200.7
$ws.Range("I85").Value = 5184
$ws.Range("J85").Value = 55267.5
$ws.Range("K85").Value = 5184
$ws.Range("L85").Value = 55267.5
$ws.Range("M85").Value = -3858
$ws.Range("N85").Value = -57919.5

$ws.Range("H107").Value = 964.125
$ws.Range("I107").Value = 983.1667
$ws.Range("K107").Value = 983.1667
$ws.Range("M107").Value = 936.8333

$ws.Range("H134").Value = 3738.7273
$ws.Range("I134").Value = 4494.2354
$ws.Range("J134").Value = 1170
$ws.Range("K134").Value = 13482.7062
$ws.Range("L134").Value = 3510
$ws.Range("M134").Value = -10947.7062
$ws.Range("N134").Value = -8580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 88
$ws.Range("I7").Value = 88
$ws.Range("K7").Value = 88
$ws.Range("M7").Value = 25

$ws.Range("H31").Value = 11787.914
$ws.Range("I31").Value = 30681.455
$ws.Range("J31").Value = 3128.375
$ws.Range("K31").Value = 30681.455
$ws.Range("L31").Value = 3128.375
$ws.Range("M31").Value = -30386.455
$ws.Range("N31").Value = -3718.375

$ws.Range("H34").Value = 11787.914
$ws.Range("I34").Value = 30681.455
$ws.Range("J34").Value = 3128.375
$ws.Range("K34").Value = 30681.455
$ws.Range("L34").Value = 3128.375
$ws.Range("M34").Value = -30479.455
$ws.Range("N34").Value = -3532.375

$ws.Range("H99").Value = 6264.2856
$ws.Range("I99").Value = 4450
$ws.Range("J99").Value = 7625
$ws.Range("K99").Value = 4450
$ws.Range("L99").Value = 7625
$ws.Range("M99").Value = -2952
$ws.Range("N99").Value = -10621

$ws.Range("H126").Value = 6264.2856
$ws.Range("I126").Value = 4450
$ws.Range("J126").Value = 7625
$ws.Range("K126").Value = 13350
$ws.Range("L126").Value = 22875
$ws.Range("M126").Value = -10880
$ws.Range("N126").Value = -27815

$ws.Range("H134").Value = 6362.5264
$ws.Range("I134").Value = 943
$ws.Range("J134").Value = 35266.668
$ws.Range("K134").Value = 2829
$ws.Range("L134").Value = 105800.004
$ws.Range("M134").Value = -294
$ws.Range("N134").Value = -110870.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4285785.5
$ws.Range("J4").Value = 7500050
$ws.Range("L4").Value = 22500150
$ws.Range("N4").Value = -22500374

$ws.Range("H5").Value = 610.4314000000001
$ws.Range("I5").Value = 533.125
$ws.Range("J5").Value = 645.7714
$ws.Range("K5").Value = 1599.375
$ws.Range("L5").Value = 1937.3142
$ws.Range("M5").Value = -1487.375
$ws.Range("N5").Value = -2161.3142

$ws.Range("H36").Value = 91094.39999999999
$ws.Range("I36").Value = 667.6
$ws.Range("J36").Value = 181521.2
$ws.Range("K36").Value = 2002.8
$ws.Range("L36").Value = 544563.6000000001
$ws.Range("M36").Value = -1833.8
$ws.Range("N36").Value = -544901.6000000001

$ws.Range("H68").Value = 5743.364
$ws.Range("J68").Value = 6545
$ws.Range("L68").Value = 19635
$ws.Range("N68").Value = -21257

$ws.Range("H71").Value = 5743.364
$ws.Range("J71").Value = 6545
$ws.Range("L71").Value = 58905
$ws.Range("N71").Value = -67017

$ws.Range("H92").Value = 31258050
$ws.Range("I92").Value = 62500350
$ws.Range("J92").Value = 15750
$ws.Range("K92").Value = 187501050
$ws.Range("L92").Value = 47250
$ws.Range("M92").Value = -187499802
$ws.Range("N92").Value = -49746

$ws.Range("H97").Value = 946.55554
$ws.Range("I97").Value = 430
$ws.Range("J97").Value = 1359.8
$ws.Range("K97").Value = 1290
$ws.Range("L97").Value = 4079.4
$ws.Range("M97").Value = -794
$ws.Range("N97").Value = -5071.4

$ws.Range("H121").Value = 4146.033
$ws.Range("I121").Value = 590
$ws.Range("J121").Value = 4857.24
$ws.Range("K121").Value = 1770
$ws.Range("L121").Value = 14571.72
$ws.Range("M121").Value = -460
$ws.Range("N121").Value = -17191.72

$ws.Range("H129").Value = 251215.8
$ws.Range("J129").Value = 334831.88
$ws.Range("L129").Value = 1004495.64
$ws.Range("N129").Value = -1014495.64

$ws.Range("H131").Value = 811.16
$ws.Range("J131").Value = 824.125
$ws.Range("L131").Value = 2472.375
$ws.Range("N131").Value = -12552.375

$ws.Range("H135").Value = 610.4314000000001
$ws.Range("I135").Value = 533.125
$ws.Range("J135").Value = 645.7714
$ws.Range("K135").Value = 4798.125
$ws.Range("L135").Value = 5811.942599999999
$ws.Range("M135").Value = -2263.125
$ws.Range("N135").Value = -10881.9426

$ws.Range("H140").Value = 1663.9524
$ws.Range("I140").Value = 1334.1666
$ws.Range("K140").Value = 4002.4998
$ws.Range("M140").Value = 1177.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3246.2
$ws.Range("I122").Value = 3056.75
$ws.Range("J122").Value = 4004
$ws.Range("K122").Value = 9170.25
$ws.Range("L122").Value = 12012
$ws.Range("M122").Value = -6720.25
$ws.Range("N122").Value = -16912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 127567.89
$ws.Range("I40").Value = 161730.14
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 161730.14
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -161594.14
$ws.Range("N40").Value = -8272

$ws.Range("H55").Value = 282.25
$ws.Range("I55").Value = 182
$ws.Range("J55").Value = 382.5
$ws.Range("K55").Value = 182
$ws.Range("L55").Value = 382.5
$ws.Range("M55").Value = -9
$ws.Range("N55").Value = -728.5

$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21996

$ws.Range("M77").ClearContents()
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69984

$ws.Range("H122").Value = 3570.5625
$ws.Range("I122").Value = 2990.5
$ws.Range("J122").Value = 4150.625
$ws.Range("K122").Value = 8971.5
$ws.Range("L122").Value = 12451.875
$ws.Range("M122").Value = -6521.5
$ws.Range("N122").Value = -17351.875

$ws.Range("H136").Value = 2624.6875
$ws.Range("I136").Value = 1799.1
$ws.Range("J136").Value = 4000.6667
$ws.Range("K136").Value = 5397.299999999999
$ws.Range("L136").Value = 12002.0001
$ws.Range("M136").Value = -2847.299999999999
$ws.Range("N136").Value = -17102.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6500
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373

$ws.Range("H122").Value = 2076.4614
$ws.Range("I122").Value = 1868.1818
$ws.Range("J122").Value = 3222
$ws.Range("K122").Value = 5604.5454
$ws.Range("L122").Value = 9666
$ws.Range("M122").Value = -3154.5454
$ws.Range("N122").Value = -14566
